$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-08 21:48:43"
$ws.Range("I2").Value = "5.9 mm"
$ws.Range("E3").Value = "2026-02-08 21:48:46"
$ws.Range("E4").Value = "2026-02-08 21:48:49"
$ws.Range("J4").Value = "1002.6 hPa"
$ws.Range("E5").Value = "2026-02-08 21:48:52"
$ws.Range("E6").Value = "2026-02-08 21:48:54"
$ws.Range("J6").Value = "1002.4 hPa"
$ws.Range("O6").Value = "10.0 °C"
$ws.Range("E7").Value = "2026-02-08 21:48:57"
$ws.Range("J7").Value = "1002.8 hPa"
$ws.Range("E8").Value = "2026-02-08 21:49:00"
$ws.Range("H8").Value = "'77%"
$ws.Range("J8").Value = "1002.7 hPa"
$ws.Range("E9").Value = "2026-02-08 21:49:02"
$ws.Range("H9").Value = "'70%"
$ws.Range("E10").Value = "2026-02-08 21:49:05"
$ws.Range("E11").Value = "2026-02-08 21:49:08"
$ws.Range("H11").Value = "'81%"
$ws.Range("O11").Value = "4.6 °C"
$ws.Range("E12").Value = "2026-02-08 21:49:11"
$ws.Range("O12").Value = "10.6 °C"
$ws.Range("E13").Value = "2026-02-08 21:49:13"
$ws.Range("J13").Value = "1004.0 hPa"
$ws.Range("E14").Value = "2026-02-08 21:49:16"
$ws.Range("K14").Value = "11.8 MJ/m2"
$ws.Range("E15").Value = "2026-02-08 21:49:19"
$ws.Range("E16").Value = "2026-02-08 21:49:21"
$ws.Range("E17").Value = "2026-02-08 21:49:24"
$ws.Range("E18").Value = "2026-02-08 21:49:26"
$ws.Range("J18").Value = "1002.8 hPa"
$ws.Range("E19").Value = "2026-02-08 21:49:29"
$ws.Range("E20").Value = "2026-02-08 21:49:32"
$ws.Range("I20").Value = "9.8 mm"
$ws.Range("E21").Value = "2026-02-08 21:49:35"
$ws.Range("J21").Value = "1003.4 hPa"
$ws.Range("K21").Value = "11.5 MJ/m2"
$ws.Range("E22").Value = "2026-02-08 21:49:38"
$ws.Range("E23").Value = "2026-02-08 21:49:40"
$ws.Range("I23").Value = "5.9 mm"
$ws.Range("E24").Value = "2026-02-08 21:49:43"
$ws.Range("H24").Value = "'80%"
$ws.Range("J24").Value = "1004.0 hPa"
$ws.Range("E25").Value = "2026-02-08 21:49:46"
$ws.Range("H25").Value = "'80%"
$ws.Range("E26").Value = "2026-02-08 21:49:49"
$ws.Range("H26").Value = "'71%"
$ws.Range("J26").Value = "1001.9 hPa"
$ws.Range("O26").Value = "3.5 °C"
$ws.Range("E27").Value = "2026-02-08 21:49:51"
$ws.Range("E28").Value = "2026-02-08 21:49:54"
$ws.Range("J28").Value = "1002.4 hPa"
$ws.Range("E29").Value = "2026-02-08 21:49:57"
$ws.Range("E30").Value = "2026-02-08 21:50:00"
$ws.Range("J30").Value = "1002.8 hPa"
$ws.Range("E31").Value = "2026-02-08 21:50:03"
$ws.Range("H31").Value = "'77%"
$ws.Range("I31").Value = "0.8 mm"
$ws.Range("J31").Value = "1002.0 hPa"
$ws.Range("N31").Value = "7.4 °C 21:28 TU"
$ws.Range("E32").Value = "2026-02-08 21:50:05"
$ws.Range("E33").Value = "2026-02-08 21:50:08"
$ws.Range("J33").Value = "1003.6 hPa"
$ws.Range("E34").Value = "2026-02-08 21:50:11"
$ws.Range("H34").Value = "'72%"
$ws.Range("E35").Value = "2026-02-08 21:50:14"
$ws.Range("J35").Value = "1004.9 hPa"
$ws.Range("E36").Value = "2026-02-08 21:50:16"
$ws.Range("J36").Value = "1002.8 hPa"
$ws.Range("E37").Value = "2026-02-08 21:50:19"
$ws.Range("J37").Value = "1003.6 hPa"
$ws.Range("E38").Value = "2026-02-08 21:50:22"
$ws.Range("O38").Value = "9.4 °C"
$ws.Range("E39").Value = "2026-02-08 21:50:24"
$ws.Range("I39").Value = "1.1 mm"
$ws.Range("E40").Value = "2026-02-08 21:50:27"
$ws.Range("H40").Value = "'84%"
$ws.Range("J40").Value = "1004.1 hPa"
$ws.Range("E41").Value = "2026-02-08 21:50:30"
$ws.Range("J41").Value = "1002.8 hPa"
$ws.Range("E42").Value = "2026-02-08 21:50:32"
$ws.Range("E43").Value = "2026-02-08 21:50:35"
$ws.Range("E44").Value = "2026-02-08 21:50:37"
$ws.Range("E45").Value = "2026-02-08 21:50:40"
$ws.Range("H45").Value = "'79%"
$ws.Range("J45").Value = "1005.0 hPa"
$ws.Range("E46").Value = "2026-02-08 21:50:43"
$ws.Range("J46").Value = "1004.5 hPa"
